$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("changeCurrent")

# --- New text labels describing the propagated-uncertainty calculation ---
# (entered in this exact order so the shared-string table comes out the same)
$ws.Range("I1").Value = "Uncertainty in V due to uncertainty in Theta:"
$ws.Range("J2").Value = "6.2845 V/rad"
$ws.Range("I2").Value = "dV/dTheta = -A ="
$ws.Range("I4").Value = "(dV/dTheta)*dTheta = "
$ws.Range("I6").Value = "Total uncertainty in V is the sum of the computer's reported uncertainty and the propogated uncertainty from the uncertainty in angle"
$ws.Range("L1").Value = "dV actual (V)"
$ws.Range("I3").Value = "dTheta = 0.5deg = "
$ws.Range("J3").Value = "0.0087 rad"
$ws.Range("J4").Value = "0.05V"

# Merged explanatory note, wrapped + centered
$ws.Range("I6:J11").Merge()
$ws.Range("I6:J11").HorizontalAlignment = -4108
$ws.Range("I6:J11").VerticalAlignment = -4108
$ws.Range("I6:J11").WrapText = $true

# Column I is wide enough to hold the labels
$ws.Columns.Item(9).ColumnWidth = 24.1640625

# --- New "actual" voltage column, offsetting the measured dV by the propagated error ---
$ws.Range("L2:L14").Formula = "=E2+0.05"
$ws.Range("L2:L14").NumberFormat = "0.00"

$ws.Range("L16").Select()
